$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Shift rows 11-21 ("Description".."Count") down to rows 12-22 to make room
# for a new "Jurisdiction" row. Work from the bottom row up so that source
# data is not clobbered before it has been copied, and clear the
# destination first since copying a blank cell does not blank the target.
for ($r = 21; $r -ge 11; $r--) {
    $src = $ws.Range("A" + $r + ":B" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":B" + ($r + 1))
    $dst.ClearContents()
    $src.Copy($dst)
}

# Fill in the new row 11 with the Jurisdiction property/value pair
$ws.Cells.Item(11, 1).Value2 = "Jurisdiction"
$ws.Cells.Item(11, 2).Value2 = "iso:code:3166:FR"

# Bump the published version number (row 3, column B)
$ws.Cells.Item(3, 2).Value2 = "0.2.0"

# Update the publication date/time (row 8, column B)
$ws.Cells.Item(8, 2).Value2 = "2023-10-20T08:59:58+00:00"
